$d = $word.ActiveDocument

# Locate the "Tree Model:" heading paragraph.
$rng = $d.Content
$rng.Find.Execute("Tree Model:", $false, $false, $false, $false, $false,
                   $true, 1, $false, $null, 0) | Out-Null

# Use the whole paragraph (including its paragraph mark) so the bold
# formatting is recorded both on the run and on the paragraph mark's
# run properties (pPr/rPr) -- matching how the other bold section
# headings ("Tunnel Texture", "Overpass Texture:", etc.) are formatted
# elsewhere in this document.
$para = $rng.Paragraphs(1)
$paraRange = $para.Range
$paraRange.Font.Bold = 1
$paraRange.Font.BoldBi = 1
